# This workbook ("Estado de Cuenta") lists, for a given NIT, one row per
# mora period (2505, 2506, ...) for each worker. A new period (2507) is
# being added on top of the existing periods, the older periods are kept
# but "pushed down" one row, and the summary totals (total mora value and
# period count) are updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for the new period by inserting a blank row right above
#    the second data row (old row 17, period 2506). This pushes the old
#    row 17 (2506) down to row 18, and everything below (rows 22/23 with
#    the signature block) down by one row as well, exactly like using
#    the worksheet's "Insert Row" command.
# ---------------------------------------------------------------------
$ws.Rows.Item(17).Insert()

# Copy the formatting of the first data row (row 16) onto the newly
# inserted row 17 so it keeps the same borders/fonts/number formats used
# by the other data rows (instead of Excel's generic blank-row format).
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Fill in the data rows. The newest period (2507) now occupies the
#    original top row (16), the period that used to be on top (2506)
#    now sits in the freshly inserted row (17), and the oldest period
#    (2505) keeps being the last data row, now at row 18.
# ---------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "9023177"
$ws.Range("D16").Value = "HUMBERTO MANUEL SANCHEZ TARRIBA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "9023177"
$ws.Range("D17").Value = "HUMBERTO MANUEL SANCHEZ TARRIBA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "9023177"
$ws.Range("D18").Value = "HUMBERTO MANUEL SANCHEZ TARRIBA"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# ---------------------------------------------------------------------
# 3. Update the summary figures at the top: the total overdue amount now
#    covers three periods instead of two (56940 * 3 = 170820), and the
#    period counter goes from 2 to 3.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3

$excel.CutCopyMode = $false
